$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "joint-datasets_after_task1"

$ws.Range("C3").Value = 1.646377249153293
$ws.Range("D3").Value = 1.336448287963867
$ws.Range("E3").Value = 51.84
$ws.Range("I3").Value = 0.0067015340924263
$ws.Range("J3").Value = 51.84

$ws.Range("C4").Value = 1.230382806670909
$ws.Range("D4").Value = 1.191968059539795
$ws.Range("E4").Value = 56.37
$ws.Range("I4").Value = 0.006039136332273484
$ws.Range("J4").Value = 56.37

$ws.Range("C5").Value = 1.045593759235071
$ws.Range("D5").Value = 0.9977485537528992
$ws.Range("E5").Value = 63.8
$ws.Range("I5").Value = 0.005048685157299042
$ws.Range("J5").Value = 63.8

$ws.Range("C6").Value = 0.9086900042027843
$ws.Range("D6").Value = 0.9067391633987427
$ws.Range("E6").Value = 67.61
$ws.Range("I6").Value = 0.004658758252859115
$ws.Range("J6").Value = 67.61

$ws.Range("C7").Value = 0.8091562198132884
$ws.Range("D7").Value = 0.8111109852790832
$ws.Range("E7").Value = 70.39
$ws.Range("I7").Value = 0.004247762632369995
$ws.Range("J7").Value = 70.39

$ws.Range("C8").Value = 0.7292847212480039
$ws.Range("D8").Value = 0.7599637389183045
$ws.Range("E8").Value = 71.42
$ws.Range("I8").Value = 0.004059910750389099
$ws.Range("J8").Value = 71.42

$ws.Range("C9").Value = 0.6632278554293574
$ws.Range("D9").Value = 0.7222894310951233
$ws.Range("E9").Value = 72.84999999999999
$ws.Range("I9").Value = 0.003899666100740433
$ws.Range("J9").Value = 72.84999999999999

$ws.Range("C10").Value = 0.5913400566091342
$ws.Range("D10").Value = 0.7514487981796265
$ws.Range("E10").Value = 73.27
$ws.Range("I10").Value = 0.003873224526643753
$ws.Range("J10").Value = 73.27

$ws.Range("C11").Value = 0.5300374897158876
$ws.Range("D11").Value = 0.7303949117660522
$ws.Range("E11").Value = 74.31
$ws.Range("I11").Value = 0.003765061503648758
$ws.Range("J11").Value = 74.31

$ws.Range("C12").Value = 0.4679656293927407
$ws.Range("D12").Value = 0.7443428874015808
$ws.Range("E12").Value = 73.97
$ws.Range("I12").Value = 0.003889581054449082
$ws.Range("J12").Value = 73.97

$ws.Range("C13").Value = 0.4076464156107027
$ws.Range("D13").Value = 0.7361296892166138
$ws.Range("E13").Value = 73.93000000000001
$ws.Range("I13").Value = 0.003958344069123268
$ws.Range("J13").Value = 73.93000000000001

$ws.Range("C14").Value = 0.3530063771471685
$ws.Range("D14").Value = 0.7372122526168823
$ws.Range("E14").Value = 75.03
$ws.Range("I14").Value = 0.003943532311916352
$ws.Range("J14").Value = 75.03

$ws.Range("C15").Value = 0.5240375629493168
$ws.Range("D15").Value = 0.7147934436798096
$ws.Range("E15").Value = 74.59999999999999
$ws.Range("I15").Value = 0.003752351200580597
$ws.Range("J15").Value = 74.59999999999999

$ws.Range("C16").Value = 0.4961007103628042
$ws.Range("D16").Value = 0.6973444700241089
$ws.Range("E16").Value = 74.84
$ws.Range("I16").Value = 0.003705506980419159
$ws.Range("J16").Value = 74.84

$ws.Range("C17").Value = 0.4750331670654063
$ws.Range("D17").Value = 0.6923460006713867
$ws.Range("E17").Value = 75.04000000000001
$ws.Range("I17").Value = 0.003689078244566917
$ws.Range("J17").Value = 75.04000000000001

$ws.Range("C18").Value = 0.45549221075311
$ws.Range("D18").Value = 0.7082119703292846
$ws.Range("E18").Value = 74.69
$ws.Range("I18").Value = 0.00375215744972229
$ws.Range("J18").Value = 74.69

$ws.Range("C19").Value = 0.4383750628451912
$ws.Range("D19").Value = 0.7176595449447631
$ws.Range("E19").Value = 74.8
$ws.Range("I19").Value = 0.003763783776760101
$ws.Range("J19").Value = 74.8

$ws.Range("C20").Value = 0.4195258226929879
$ws.Range("D20").Value = 0.7028288602828979
$ws.Range("E20").Value = 75.36
$ws.Range("I20").Value = 0.003722845929861069
$ws.Range("J20").Value = 75.36

$ws.Range("C21").Value = 0.3996751639307762
$ws.Range("D21").Value = 0.7138570189476013
$ws.Range("E21").Value = 75.20999999999999
$ws.Range("I21").Value = 0.00374254395365715
$ws.Range("J21").Value = 75.20999999999999

$ws.Range("C22").Value = 0.3836170034140957
$ws.Range("D22").Value = 0.714115834236145
$ws.Range("E22").Value = 75.3
$ws.Range("I22").Value = 0.003790253400802612
$ws.Range("J22").Value = 75.3

$ws.Range("C23").Value = 0.4398005764095151
$ws.Range("D23").Value = 0.6904516339302063
$ws.Range("E23").Value = 75.31999999999999
$ws.Range("I23").Value = 0.003658279645442962
$ws.Range("J23").Value = 75.31999999999999

$ws.Range("C24").Value = 0.4335430353271718
$ws.Range("D24").Value = 0.6886540412902832
$ws.Range("E24").Value = 75.42
$ws.Range("I24").Value = 0.003655935353040695
$ws.Range("J24").Value = 75.42

$ws.Range("C25").Value = 0.4283898581047447
$ws.Range("D25").Value = 0.6951533436775208
$ws.Range("E25").Value = 75.38
$ws.Range("I25").Value = 0.003668103164434433
$ws.Range("J25").Value = 75.38

$ws.Range("C26").Value = 0.424234146366314
$ws.Range("D26").Value = 0.6972607731819153
$ws.Range("E26").Value = 75.53
$ws.Range("I26").Value = 0.0036801553606987
$ws.Range("J26").Value = 75.53

$ws.Range("C27").Value = 0.4199508109871222
$ws.Range("D27").Value = 0.6900336027145386
$ws.Range("E27").Value = 75.5
$ws.Range("I27").Value = 0.003668453359603882
$ws.Range("J27").Value = 75.5

$ws.Range("C28").Value = 0.4158378115722111
$ws.Range("D28").Value = 0.6927130937576294
$ws.Range("E28").Value = 75.28
$ws.Range("I28").Value = 0.003677964514493942
$ws.Range("J28").Value = 75.28

$ws.Range("C29").Value = 0.4101625398105505
$ws.Range("D29").Value = 0.6903794884681702
$ws.Range("E29").Value = 75.53
$ws.Range("I29").Value = 0.003657609480619431
$ws.Range("J29").Value = 75.53

$ws.Range("C30").Value = 0.425295884755193
$ws.Range("D30").Value = 0.6905232787132263
$ws.Range("E30").Value = 75.47
$ws.Range("I30").Value = 0.003653058338165283
$ws.Range("J30").Value = 75.47

$ws.Range("C31").Value = 0.4238404433338009
$ws.Range("D31").Value = 0.6902877926826477
$ws.Range("E31").Value = 75.3
$ws.Range("I31").Value = 0.003653685906529426
$ws.Range("J31").Value = 75.3

$ws.Range("C32").Value = 0.4227969134340481
$ws.Range("D32").Value = 0.6905940175056458
$ws.Range("E32").Value = 75.3
$ws.Range("I32").Value = 0.003654783833026886
$ws.Range("J32").Value = 75.3

$ws.Range("C33").Value = 0.4217399243189364
$ws.Range("D33").Value = 0.6904921770095825
$ws.Range("E33").Value = 75.34
$ws.Range("I33").Value = 0.003655805158615112
$ws.Range("J33").Value = 75.34

$ws.Range("C34").Value = 0.4207482026547802
$ws.Range("D34").Value = 0.6910329818725586
$ws.Range("E34").Value = 75.45999999999999
$ws.Range("I34").Value = 0.003657100480794906
$ws.Range("J34").Value = 75.45999999999999

$ws.Range("C35").Value = 0.424675516206391
$ws.Range("D35").Value = 0.6897433876991272
$ws.Range("E35").Value = 75.40000000000001
$ws.Range("I35").Value = 0.00365116440653801
$ws.Range("J35").Value = 75.40000000000001

$ws.Range("C36").Value = 0.4240329271676589
$ws.Range("D36").Value = 0.689545464515686
$ws.Range("E36").Value = 75.33
$ws.Range("I36").Value = 0.003651230055093765
$ws.Range("J36").Value = 75.33

$ws.Range("C37").Value = 0.4237239598011484
$ws.Range("D37").Value = 0.68962721824646
$ws.Range("E37").Value = 75.31999999999999
$ws.Range("I37").Value = 0.00365224147439003
$ws.Range("J37").Value = 75.31999999999999

$ws.Range("C38").Value = 0.4234603730999694
$ws.Range("D38").Value = 0.6898393869400025
$ws.Range("E38").Value = 75.37
$ws.Range("I38").Value = 0.003653755432367325
$ws.Range("J38").Value = 75.37

$ws.Range("C39").Value = 0.4232553572070842
$ws.Range("D39").Value = 0.6899454116821289
$ws.Range("E39").Value = 75.39
$ws.Range("I39").Value = 0.00365321980714798
$ws.Range("J39").Value = 75.39
